$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 725.0417
$ws.Range("I28").Value = 403.94116
$ws.Range("J28").Value = 1504.8572
$ws.Range("K28").Value = 403.94116
$ws.Range("L28").Value = 1504.8572
$ws.Range("M28").Value = 81.05883999999998
$ws.Range("N28").Value = -2474.8572

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 250
$ws.Range("I29").Value = 140
$ws.Range("J29").Value = 800
$ws.Range("K29").Value = 420
$ws.Range("L29").Value = 2400
$ws.Range("M29").Value = -139
$ws.Range("N29").Value = -2962

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 979.9091
$ws.Range("I102").Value = 957.9
$ws.Range("J102").Value = 1200
$ws.Range("K102").Value = 957.9
$ws.Range("L102").Value = 1200
$ws.Range("M102").Value = 664.1
$ws.Range("N102").Value = -4444

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 6503.4814
$ws.Range("I122").Value = 8119.7
$ws.Range("J122").Value = 1885.7142
$ws.Range("K122").Value = 24359.1
$ws.Range("L122").Value = 5657.142599999999
$ws.Range("M122").Value = -21909.1
$ws.Range("N122").Value = -10557.1426

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 9618134
$ws.Range("I132").Value = 15627347
$ws.Range("J132").Value = 3392.6
$ws.Range("K132").Value = 46882041
$ws.Range("L132").Value = 10177.8
$ws.Range("M132").Value = -46879511
$ws.Range("N132").Value = -15237.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4528.8887
$ws.Range("I105").Value = 4000
$ws.Range("J105").Value = 4549.231
$ws.Range("K105").Value = 4000
$ws.Range("L105").Value = 4549.231
$ws.Range("M105").Value = -2253
$ws.Range("N105").Value = -8043.231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4259.304
$ws.Range("I134").Value = 3351.0527
$ws.Range("J134").Value = 4898.4443
$ws.Range("K134").Value = 10053.1581
$ws.Range("L134").Value = 14695.3329
$ws.Range("M134").Value = -7518.158100000001
$ws.Range("N134").Value = -19765.3329

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2230.13
$ws.Range("I31").Value = 1920.561
$ws.Range("J31").Value = 2582.6943
$ws.Range("K31").Value = 1920.561
$ws.Range("L31").Value = 2582.6943
$ws.Range("M31").Value = -1625.561
$ws.Range("N31").Value = -3172.6943

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2230.13
$ws.Range("I34").Value = 1920.561
$ws.Range("J34").Value = 2582.6943
$ws.Range("K34").Value = 1920.561
$ws.Range("L34").Value = 2582.6943
$ws.Range("M34").Value = -1718.561
$ws.Range("N34").Value = -2986.6943

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1826.5652
$ws.Range("I99").Value = 1515.5385
$ws.Range("J99").Value = 2230.9
$ws.Range("K99").Value = 1515.5385
$ws.Range("L99").Value = 2230.9
$ws.Range("M99").Value = -17.53850000000011
$ws.Range("N99").Value = -5226.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1826.5652
$ws.Range("I126").Value = 1515.5385
$ws.Range("J126").Value = 2230.9
$ws.Range("K126").Value = 4546.6155
$ws.Range("L126").Value = 6692.700000000001
$ws.Range("M126").Value = -2076.6155
$ws.Range("N126").Value = -11632.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1402739.5
$ws.Range("I134").Value = 2423.1333
$ws.Range("J134").Value = 11905112
$ws.Range("K134").Value = 7269.3999
$ws.Range("L134").Value = 35715336
$ws.Range("M134").Value = -4734.3999
$ws.Range("N134").Value = -35720406

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 55556350
$ws.Range("I113").Value = 166667230
$ws.Range("J113").Value = 914.8333
$ws.Range("K113").Value = 500001690
$ws.Range("L113").Value = 2744.4999
$ws.Range("M113").Value = -499999520
$ws.Range("N113").Value = -7084.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 511.7857
$ws.Range("I121").Value = 484
$ws.Range("J121").Value = 581.25
$ws.Range("K121").Value = 1452
$ws.Range("L121").Value = 1743.75
$ws.Range("M121").Value = -142
$ws.Range("N121").Value = -4363.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3952.9524
$ws.Range("I102").Value = 4645.6
$ws.Range("J102").Value = 2221.3333
$ws.Range("K102").Value = 4645.6
$ws.Range("L102").Value = 2221.3333
$ws.Range("M102").Value = -3023.6
$ws.Range("N102").Value = -5465.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2803.2
$ws.Range("I122").Value = 2724.8
$ws.Range("J122").Value = 2960
$ws.Range("K122").Value = 8174.400000000001
$ws.Range("L122").Value = 8880
$ws.Range("M122").Value = -5724.400000000001
$ws.Range("N122").Value = -13780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3599.647
$ws.Range("I126").Value = 2324.875
$ws.Range("J126").Value = 4732.778
$ws.Range("K126").Value = 6974.625
$ws.Range("L126").Value = 14198.334
$ws.Range("M126").Value = -4504.625
$ws.Range("N126").Value = -19138.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7394.1
$ws.Range("I40").Value = 11276
$ws.Range("J40").Value = 3512.2
$ws.Range("K40").Value = 11276
$ws.Range("L40").Value = 3512.2
$ws.Range("M40").Value = -11140
$ws.Range("N40").Value = -3784.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1600.7646
$ws.Range("I68").Value = 1886.25
$ws.Range("J68").Value = 1347
$ws.Range("K68").Value = 1886.25
$ws.Range("L68").Value = 1347
$ws.Range("M68").Value = -1137.25
$ws.Range("N68").Value = -2845

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1600.7646
$ws.Range("I71").Value = 1886.25
$ws.Range("J71").Value = 1347
$ws.Range("K71").Value = 9431.25
$ws.Range("L71").Value = 6735
$ws.Range("M71").Value = -5687.25
$ws.Range("N71").Value = -14223

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1413.3
$ws.Range("I82").Value = 1280.3529
$ws.Range("J82").Value = 2166.6667
$ws.Range("K82").Value = 1280.3529
$ws.Range("L82").Value = 2166.6667
$ws.Range("M82").Value = -919.3529000000001
$ws.Range("N82").Value = -2888.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1413.3
$ws.Range("I85").Value = 1280.3529
$ws.Range("J85").Value = 2166.6667
$ws.Range("K85").Value = 1280.3529
$ws.Range("L85").Value = 2166.6667
$ws.Range("M85").Value = -32.35290000000009
$ws.Range("N85").Value = -4662.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8935504
$ws.Range("I132").Value = 4588.095
$ws.Range("J132").Value = 35728252
$ws.Range("K132").Value = 13764.285
$ws.Range("L132").Value = 107184756
$ws.Range("M132").Value = -11234.285
$ws.Range("N132").Value = -107189816

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1894.8518
$ws.Range("I126").Value = 1577.125
$ws.Range("J126").Value = 4436.6665
$ws.Range("K126").Value = 4731.375
$ws.Range("L126").Value = 13309.9995
$ws.Range("M126").Value = -2261.375
$ws.Range("N126").Value = -18249.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1414.1731
$ws.Range("I132").Value = 1170.66
$ws.Range("K132").Value = 3511.98
$ws.Range("M132").Value = -981.9800000000005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1469.1957
$ws.Range("I136").Value = 888
$ws.Range("J136").Value = 2944.5386
$ws.Range("K136").Value = 2664
$ws.Range("L136").Value = 8833.6158
$ws.Range("M136").Value = -114
$ws.Range("N136").Value = -13933.6158
